$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 17-20 which are no longer part of the dataset
$ws.Range("A17:A20").EntireRow.Delete()

# Update B/C values for rows 2-16
$values = @(
    @(0, 0, 0),
    @(1, 4.738603388729921, 0.001156177677677678),
    @(2, 7.339947914712358, 0.001849884284284284),
    @(3, 11.23519217776038, 0.003006061961961962),
    @(4, 16.04563476924854, 0.004624710710710711),
    @(5, 19.54107307980702, 0.006012123923923924),
    @(6, 23.00569908379239, 0.007630772672672673),
    @(7, 27.86900011229479, 0.01063683463463463),
    @(8, 30.79460126660187, 0.01318042552552552),
    @(9, 34.3398220222001, 0.01780513623623623),
    @(10, 36.39422832549062, 0.02196737587587588),
    @(11, 37.79655091322297, 0.02589837997997998),
    @(12, 39.00197108910267, 0.0312167972972973),
    @(13, 39.62802055302561, 0.03815386336336336),
    @(14, 39.66155870212616, 0.042316103003003)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
    $ws.Cells.Item($row, 3).Value = $values[$i][2]
}
